$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3430
$ws.Range("I19").Value = 1992
$ws.Range("J19").Value = 3909.3333
$ws.Range("K19").Value = 1992
$ws.Range("L19").Value = 3909.3333
$ws.Range("M19").Value = -1817
$ws.Range("N19").Value = -4259.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4422.8823
$ws.Range("I51").Value = 3799.4614
$ws.Range("J51").Value = 6449
$ws.Range("K51").Value = 3799.4614
$ws.Range("L51").Value = 6449
$ws.Range("M51").Value = -3315.4614
$ws.Range("N51").Value = -7417

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3599.6667
$ws.Range("I80").Value = 3514.5715
$ws.Range("J80").Value = 3674.125
$ws.Range("K80").Value = 10543.7145
$ws.Range("L80").Value = 11022.375
$ws.Range("M80").Value = -9545.7145
$ws.Range("N80").Value = -13018.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3599.6667
$ws.Range("I83").Value = 3514.5715
$ws.Range("J83").Value = 3674.125
$ws.Range("K83").Value = 31631.1435
$ws.Range("L83").Value = 33067.125
$ws.Range("M83").Value = -26639.1435
$ws.Range("N83").Value = -43051.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 20001638
$ws.Range("I101").Value = 25000798
$ws.Range("K101").Value = 75002394
$ws.Range("M101").Value = -75000772

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1332.6666
$ws.Range("I132").Value = 1349.5
$ws.Range("K132").Value = 4048.5
$ws.Range("M132").Value = -1518.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3888.027
$ws.Range("I138").Value = 1883.0834
$ws.Range("J138").Value = 4850.4
$ws.Range("K138").Value = 5649.2502
$ws.Range("L138").Value = 14551.2
$ws.Range("M138").Value = -509.2502000000004
$ws.Range("N138").Value = -24831.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2336.125
$ws.Range("I141").Value = 2198.4
$ws.Range("J141").Value = 2565.6667
$ws.Range("K141").Value = 6595.200000000001
$ws.Range("L141").Value = 7697.000100000001
$ws.Range("M141").Value = -1415.200000000001
$ws.Range("N141").Value = -18057.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2011
$ws.Range("I2").Value = 2011
$ws.Range("K2").Value = 2011
$ws.Range("M2").Value = -1898

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 334666.34
$ws.Range("I45").Value = 334666.34
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 334666.34
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -334289.34
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6358.5454
$ws.Range("J63").Value = 6397
$ws.Range("L63").Value = 6397
$ws.Range("N63").Value = -7769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6358.5454
$ws.Range("J66").Value = 6397
$ws.Range("L66").Value = 31985
$ws.Range("N66").Value = -38849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1791.25
$ws.Range("I97").Value = 1755.5
$ws.Range("J97").Value = 1898.5
$ws.Range("K97").Value = 1755.5
$ws.Range("L97").Value = 1898.5
$ws.Range("M97").Value = -1259.5
$ws.Range("N97").Value = -2890.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2011
$ws.Range("I116").Value = 2011
$ws.Range("K116").Value = 2011
$ws.Range("M116").Value = 283

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2011
$ws.Range("I3").Value = 2011
$ws.Range("K3").Value = 2011
$ws.Range("M3").Value = -1897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 950.875
$ws.Range("J64").Value = 956
$ws.Range("L64").Value = 956
$ws.Range("N64").Value = -1406

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 950.875
$ws.Range("J67").Value = 956
$ws.Range("L67").Value = 956
$ws.Range("N67").Value = -2516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1852.5454
$ws.Range("I86").Value = 1768.5714
$ws.Range("K86").Value = 1768.5714
$ws.Range("M86").Value = -645.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1852.5454
$ws.Range("I89").Value = 1768.5714
$ws.Range("K89").Value = 8842.857
$ws.Range("M89").Value = -3226.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3468.75
$ws.Range("I134").Value = 3468.75
$ws.Range("K134").Value = 10406.25
$ws.Range("M134").Value = -7871.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2488.7693
$ws.Range("I16").Value = 2740.2222
$ws.Range("K16").Value = 2740.2222
$ws.Range("M16").Value = -2453.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 539.1429000000001
$ws.Range("I22").Value = 522
$ws.Range("J22").Value = 562
$ws.Range("K22").Value = 522
$ws.Range("L22").Value = 562
$ws.Range("M22").Value = -172
$ws.Range("N22").Value = -1262

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 163545.86
$ws.Range("I94").Value = 282250
$ws.Range("J94").Value = 5273.6665
$ws.Range("K94").Value = 282250
$ws.Range("L94").Value = 5273.6665
$ws.Range("M94").Value = -281799
$ws.Range("N94").Value = -6175.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2488.7693
$ws.Range("I113").Value = 2740.2222
$ws.Range("K113").Value = 2740.2222
$ws.Range("M113").Value = -570.2222000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 372.1
$ws.Range("I12").Value = 420.5
$ws.Range("J12").Value = 299.5
$ws.Range("K12").Value = 1261.5
$ws.Range("L12").Value = 898.5
$ws.Range("M12").Value = -1088.5
$ws.Range("N12").Value = -1244.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1400
$ws.Range("I58").Value = 1400
$ws.Range("K58").Value = 4200
$ws.Range("M58").Value = -4072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 10636.444
$ws.Range("I70").Value = 1432.5
$ws.Range("K70").Value = 4297.5
$ws.Range("M70").Value = -3982.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 10636.444
$ws.Range("I73").Value = 1432.5
$ws.Range("K73").Value = 4297.5
$ws.Range("M73").Value = -3205.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 825
$ws.Range("I109").Value = 825
$ws.Range("K109").Value = 2475
$ws.Range("M109").Value = -1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2500360
$ws.Range("I11").Value = 3143371.2
$ws.Range("K11").Value = 3143371.2
$ws.Range("M11").Value = -3143232.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3397.3333
$ws.Range("I102").Value = 3397.3333
$ws.Range("K102").Value = 3397.3333
$ws.Range("M102").Value = -1775.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2257.9167
$ws.Range("I107").Value = 818.1667
$ws.Range("J107").Value = 3697.6667
$ws.Range("K107").Value = 818.1667
$ws.Range("L107").Value = 3697.6667
$ws.Range("M107").Value = 1101.8333
$ws.Range("N107").Value = -7537.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3388.8
$ws.Range("J126").Value = 3125
$ws.Range("L126").Value = 9375
$ws.Range("N126").Value = -14315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7416.1665
$ws.Range("I7").Value = 3334.3333
$ws.Range("J7").Value = 8776.777
$ws.Range("K7").Value = 3334.3333
$ws.Range("L7").Value = 8776.777
$ws.Range("M7").Value = -3222.3333
$ws.Range("N7").Value = -9000.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 787.05884
$ws.Range("I55").Value = 421.0909
$ws.Range("J55").Value = 1458
$ws.Range("K55").Value = 421.0909
$ws.Range("L55").Value = 1458
$ws.Range("M55").Value = -248.0909
$ws.Range("N55").Value = -1804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1995
$ws.Range("I82").Value = 700
$ws.Range("K82").Value = 700
$ws.Range("M82").Value = -339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1995
$ws.Range("I85").Value = 700
$ws.Range("K85").Value = 700
$ws.Range("M85").Value = 548

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7416.1665
$ws.Range("I126").Value = 3334.3333
$ws.Range("J126").Value = 8776.777
$ws.Range("K126").Value = 10002.9999
$ws.Range("L126").Value = 26330.331
$ws.Range("M126").Value = -7532.999899999999
$ws.Range("N126").Value = -31270.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1885.125
$ws.Range("I136").Value = 926.17645
$ws.Range("K136").Value = 2778.52935
$ws.Range("M136").Value = -228.5293500000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1125000
$ws.Range("I8").Value = 1125000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1125000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1124860
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 74110
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 74110
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3328.9
$ws.Range("I132").Value = 3476
$ws.Range("K132").Value = 10428
$ws.Range("M132").Value = -7898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 360071.66
$ws.Range("J135").Value = 360071.66
$ws.Range("L135").Value = 360071.66
$ws.Range("N135").Value = -370211.66
